$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.903.91"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.74"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.31"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06423"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.71"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.665.88"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.313"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5477"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7904"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.11"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.988.91"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.62"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.430"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.076"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.856"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1147"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.907"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.242"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05030"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.207"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.545"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8939"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.597"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5548"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.134.75"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01565"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.671"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8156"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.75"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.786.01"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4544"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05096"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09568"
$ws.Range("E51").Value = "  +3.05%  "
